# Adapt column header formatting to respective input file names.
# "<field>_old"/"<field>_new" -> "<field>_FV2304"/"<field>_FV2310"
# then (re)expose the data range as an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J carry the "_old" -> "_FV2304" suffixed headers.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $baseNames[$i] + "_FV2304"
}

# Column K is the unchanged "diff" header.
$ws.Cells.Item(1, 11).Value2 = "diff"

# Columns L-U carry the "_new" -> "_FV2310" suffixed headers.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value2 = $baseNames[$i] + "_FV2310"
}

# Turn the whole used range into an Excel Table (ListObject) now that the
# header row carries the new names.
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1).
[void]$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Cells.Item(1, 1).Select()
